$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '''27.000.49'
$ws.Range("E2").Value = '  -0.47%  '

# Row 3
$ws.Range("D3").Value = '''1.827.53'
$ws.Range("E3").Value = '  +0.22%  '

# Row 4
$ws.Range("D4").Value = '''1.007'
$ws.Range("E4").Value = '  -0.44%  '

# Row 5
$ws.Range("D5").Value = '''312.32'
$ws.Range("E5").Value = '  -0.01%  '

# Row 6
$ws.Range("D6").Value = '''1.006'

# Row 7
$ws.Range("D7").Value = '''0.4584'
$ws.Range("E7").Value = '  -0.79%  '

# Row 8
$ws.Range("D8").Value = '''0.3703'

# Row 9
$ws.Range("D9").Value = '''0.07319'
$ws.Range("E9").Value = '  +0.25%  '

# Row 10
$ws.Range("D10").Value = '''0.8730'
$ws.Range("E10").Value = '  +0.31%  '

# Row 11
$ws.Range("D11").Value = '''0.07943'
$ws.Range("E11").Value = '  +4.05%  '

# Row 12
$ws.Range("D12").Value = '''19.78'
$ws.Range("E12").Value = '  -1.55%  '

# Row 13
$ws.Range("D13").Value = '''1.897.67'
$ws.Range("E13").Value = '  +0.67%  '

# Row 14
$ws.Range("D14").Value = '''6.565'
$ws.Range("E14").Value = '  +1.50%  '

# Row 15
$ws.Range("D15").Value = '''5.330'
$ws.Range("E15").Value = '  -0.12%  '

# Row 16
$ws.Range("D16").Value = '''91.41'
$ws.Range("E16").Value = '  -1.02%  '

# Row 17
$ws.Range("E17").Value = '  -0.26%  '

# Row 18
$ws.Range("D18").Value = '''0.000008903'
$ws.Range("E18").Value = '  +3.07%  '

# Row 19
$ws.Range("E19").Value = '  -0.37%  '

# Row 20
$ws.Range("D20").Value = '''14.71'
$ws.Range("E20").Value = '  +1.69%  '

# Row 21
$ws.Range("D21").Value = '''27.196.03'
$ws.Range("E21").Value = '  -0.88%  '

# Row 22
$ws.Range("D22").Value = '''5.097'
$ws.Range("E22").Value = '  -2.26%  '

# Row 23
$ws.Range("E23").Value = '  -0.01%  '

# Row 24
$ws.Range("D24").Value = '''2.100.56'
$ws.Range("E24").Value = '  +0.27%  '

# Row 25
$ws.Range("D25").Value = '''153.02'
$ws.Range("E25").Value = '  +1.05%  '

# Row 26
$ws.Range("E26").Value = '  -1.52%  '

# Row 27
$ws.Range("D27").Value = '''18.38'
$ws.Range("E27").Value = '  +1.67%  '

# Row 28
$ws.Range("D28").Value = '''2.050'
$ws.Range("E28").Value = '  -1.57%  '

# Row 29
$ws.Range("D29").Value = '''5.136'
$ws.Range("E29").Value = '  +0.76%  '

# Row 30
$ws.Range("D30").Value = '''115.08'
$ws.Range("E30").Value = '  -0.63%  '

# Row 31
$ws.Range("D31").Value = '''0.08860'
$ws.Range("E31").Value = '  -0.44%  '

# Row 32
$ws.Range("D32").Value = '''2.963'
$ws.Range("E32").Value = '  +0.10%  '

# Row 33
$ws.Range("D33").Value = '''0.7301'
$ws.Range("E33").Value = '  -1.28%  '

# Row 34
$ws.Range("D34").Value = '''4.448'
$ws.Range("E34").Value = '  -0.01%  '

# Row 35
$ws.Range("E35").Value = '  -0.50%  '

# Row 36
$ws.Range("D36").Value = '''1.074'
$ws.Range("E36").Value = '  +0.26%  '

# Row 37
$ws.Range("B37").Value = 'RenderToken'
$ws.Range("C37").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D37").Value = '''2.446'
$ws.Range("E37").Value = '  -1.43%  '

# Row 38
$ws.Range("B38").Value = 'VeChain'
$ws.Range("C38").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D38").Value = '''0.01943'
$ws.Range("E38").Value = '  +1.72%  '

# Row 39
$ws.Range("B39").Value = 'Hedera'
$ws.Range("C39").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D39").Value = '''0.05233'
$ws.Range("E39").Value = '  -0.19%  '

# Row 40
$ws.Range("D40").Value = '''2.947'
$ws.Range("E40").Value = '  +0.68%  '

# Row 41
$ws.Range("D41").Value = '''7.143'
$ws.Range("E41").Value = '  -0.24%  '

# Row 42
$ws.Range("D42").Value = '''0.5143'
$ws.Range("E42").Value = '  -0.90%  '

# Row 43
$ws.Range("D43").Value = '''0.1630'
$ws.Range("E43").Value = '  +0.23%  '

# Row 44
$ws.Range("D44").Value = '''8.206'
$ws.Range("E44").Value = '  -0.80%  '

# Row 45
$ws.Range("D45").Value = '''0.4837'
$ws.Range("E45").Value = '  +0.15%  '

# Row 46
$ws.Range("B46").Value = 'EnergySwap'
$ws.Range("C46").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D46").Value = '''10.23'
$ws.Range("E46").Value = '  +0.56%  '

# Row 47
$ws.Range("B47").Value = 'PaxDollar'
$ws.Range("C47").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D47").Value = '''1.006'
$ws.Range("E47").Value = '  -0.42%  '

# Row 48
$ws.Range("D48").Value = '''102.49'
$ws.Range("E48").Value = '  -0.82%  '

# Row 49
$ws.Range("D49").Value = '''1.626'
$ws.Range("E49").Value = '  -0.45%  '

# Row 50
$ws.Range("D50").Value = '''0.06212'
$ws.Range("E50").Value = '  -0.90%  '

# Row 51
$ws.Range("D51").Value = '''64.82'
$ws.Range("E51").Value = '  +0.57%  '

